$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.670.42'
$ws.Range("E2").Value = '  -0.76%  '

$ws.Range("D3").Value = '3.646.13'
$ws.Range("E3").Value = '  -1.12%  '

$ws.Range("E4").Value = '  +0.41%  '

$ws.Range("D5").Value = "'586.19"
$ws.Range("E5").Value = '  -1.28%  '

$ws.Range("D6").Value = "'186.61"
$ws.Range("E6").Value = '  +2.56%  '

$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = '  -2.19%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").Value = "'0.683"
$ws.Range("E9").Value = '  -4.27%  '

$ws.Range("D10").Value = "'0.148"
$ws.Range("E10").Value = '  -8.96%  '

$ws.Range("D11").Value = "'54.94"
$ws.Range("E11").Value = '  -1.70%  '

$ws.Range("D12").Value = "'0.0000260"
$ws.Range("E12").Value = '  -10.25%  '

$ws.Range("D13").Value = "'10.03"
$ws.Range("E13").Value = '  -3.05%  '

$ws.Range("D14").Value = '4.244.97'
$ws.Range("E14").Value = '  -0.59%  '

$ws.Range("D15").Value = '3.662.39'
$ws.Range("E15").Value = '  -0.59%  '

$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '67.661.14'
$ws.Range("E17").Value = '  -0.44%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = "'18.52"
$ws.Range("E18").Value = '  -3.82%  '

$ws.Range("D19").Value = "'1.08"
$ws.Range("E19").Value = '  -3.12%  '

$ws.Range("D20").Value = "'12.33"
$ws.Range("E20").Value = '  -3.46%  '

$ws.Range("D21").Value = "'396.31"
$ws.Range("E21").Value = '  -2.91%  '

$ws.Range("D22").Value = "'4.34"
$ws.Range("E22").Value = '  -4.44%  '

$ws.Range("D23").Value = "'86.32"
$ws.Range("E23").Value = '  -2.61%  '

$ws.Range("D24").Value = "'2.89"
$ws.Range("E24").Value = '  -4.19%  '

$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").Value = "'10.62"
$ws.Range("E25").Value = '  -2.90%  '

$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = "'12.28"
$ws.Range("E26").Value = '  -3.72%  '

$ws.Range("E27").Value = '  +0.67%  '

$ws.Range("D28").Value = "'3.62"
$ws.Range("E28").Value = '  -6.06%  '

$ws.Range("D29").Value = "'9.11"
$ws.Range("E29").Value = '  -3.05%  '

$ws.Range("D30").Value = "'31.45"
$ws.Range("E30").Value = '  -3.85%  '

$ws.Range("D31").Value = "'6.94"
$ws.Range("E31").Value = '  -4.14%  '

$ws.Range("D32").Value = "'67.13"
$ws.Range("E32").Value = '  +4.32%  '

$ws.Range("D33").Value = "'12.00"
$ws.Range("E33").Value = '  -3.46%  '

$ws.Range("D34").Value = "'600.79"
$ws.Range("E34").Value = '  -0.07%  '

$ws.Range("D35").Value = "'42.86"
$ws.Range("E35").Value = '  -1.26%  '

$ws.Range("D36").Value = "'0.113"
$ws.Range("E36").Value = '  -3.12%  '

$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("E38").Value = '  +0.28%  '

$ws.Range("D39").Value = "'0.382"
$ws.Range("E39").Value = '  -4.01%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = "'0.135"
$ws.Range("E40").Value = '  -0.97%  '

$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0729'
$ws.Range("E41").Value = '  -17.41%  '

$ws.Range("D42").Value = "'2.82"
$ws.Range("E42").Value = '  -5.70%  '

$ws.Range("D43").Value = "'0.0414"
$ws.Range("E43").Value = '  -5.06%  '

$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = "'3.19"
$ws.Range("E44").Value = '  +1.26%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = "'2.44"
$ws.Range("E45").Value = '  -12.24%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = "'0.133"
$ws.Range("E46").Value = '  -0.59%  '

$ws.Range("D47").Value = '2.729.01'
$ws.Range("E47").Value = '  -0.11%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = "'140.77"
$ws.Range("E48").Value = '  -0.86%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = "'8.62"
$ws.Range("E49").Value = '  -6.16%  '

$ws.Range("E50").Value = '  -5.75%  '

$ws.Range("D51").Value = "'2.63"
$ws.Range("E51").Value = '  -4.63%  '
